$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.476.77'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '2.492.87'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''569.81'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").Value = '''166.72'
$ws.Range("E6").Value = '  +1.24%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("D11").Value = '''0.350'
$ws.Range("E11").Value = '  -1.05%  '
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '2.949.54'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '69.389.51'
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").Value = '''24.22'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '2.488.92'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").Value = '''11.21'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''352.60'
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''7.38'
$ws.Range("E20").Value = '  -3.22%  '
$ws.Range("E21").Value = '  +0.56%  '
$ws.Range("D22").Value = '''1.91'
$ws.Range("E22").Value = '  -3.65%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("D25").Value = '''3.79'
$ws.Range("E25").Value = '  -2.56%  '
$ws.Range("D26").Value = '2.621.56'
$ws.Range("E26").Value = '  -1.04%  '
$ws.Range("E27").Value = '  -2.00%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("D30").Value = '''7.58'
$ws.Range("E30").Value = '  -1.99%  '
$ws.Range("D31").Value = '''442.98'
$ws.Range("E31").Value = '  -3.13%  '
$ws.Range("D32").Value = '''1.20'
$ws.Range("E32").Value = '  -2.21%  '
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("D35").Value = '''3.12'
$ws.Range("E35").Value = '  +104.75%  '
$ws.Range("D36").Value = '''154.23'
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("E37").Value = '  -1.92%  '
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").Value = '''18.13'
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '''0.313'
$ws.Range("E41").Value = '  -0.89%  '
$ws.Range("E42").Value = '  -1.06%  '
$ws.Range("E43").Value = '  -0.66%  '
$ws.Range("D44").Value = '''2.19'
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("E45").Value = '  -3.87%  '
$ws.Range("D46").Value = '''138.67'
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").Value = '''0.505'
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("D49").Value = '''0.0723'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("D51").Value = '''0.0923'
$ws.Range("E51").Value = '  -0.43%  '
